# Daily attendance processing - 2025-12-31 07:58:35
# Reorders the "Recorded By" (column G) value lists so that email
# addresses (tokens containing "@") come first, sorted alphabetically,
# followed by the remaining non-email tokens (e.g. "System"/"system") in
# their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reorder-RecordedBy {
    param([string]$val)

    $parts = $val -split ", "
    $emails = @()
    $nonEmails = @()

    foreach ($p in $parts) {
        if ($p -like "*@*") {
            $emails += $p
        } else {
            $nonEmails += $p
        }
    }

    $emailsSorted = @($emails | Sort-Object)
    $combo = @($emailsSorted) + @($nonEmails)

    return ($combo -join ", ")
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $current = $cell.Value2
    if ($current -ne $null -and $current -ne "") {
        $new = Reorder-RecordedBy $current
        if ($new -ne $current) {
            $cell.Value2 = $new
        }
    }
}
